$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.974945333333333
$ws.Range("H2").Value = 23.924836
$ws.Range("I2").Value = 0.01228328074175765
$ws.Range("J2").Value = 0.01228328074175765
$ws.Range("M2").Value = 0.3360566666666667
$ws.Range("N2").Value = 1.00817
$ws.Range("O2").Value = 0.01570866217798777
$ws.Range("P2").Value = 0.01570866217798777
$ws.Range("Q2").Value = 2.680033545568889
$ws.Range("R2").Value = 24.12030191012
$ws.Range("S2").Value = 0.0001929539076096539
$ws.Range("T2").Value = 0.0001929539076096539
$ws.Range("G3").Value = 7.974945333333333
$ws.Range("H3").Value = 23.924836
$ws.Range("I3").Value = 0.01228328074175765
$ws.Range("J3").Value = 0.01228328074175765
$ws.Range("O3").Value = 0.109316751024163
$ws.Range("P3").Value = 0.1093167510241629
$ws.Range("Q3").Value = 18.65038260533067
$ws.Range("R3").Value = 167.853443447976
$ws.Range("S3").Value = 0.001342768342606616
$ws.Range("T3").Value = 0.001342768342606616
$ws.Range("G4").Value = 7.974945333333333
$ws.Range("H4").Value = 23.924836
$ws.Range("I4").Value = 0.01228328074175765
$ws.Range("J4").Value = 0.01228328074175765
$ws.Range("M4").Value = 18.491866
$ws.Range("N4").Value = 55.47559800000001
$ws.Range("O4").Value = 0.864385399390831
$ws.Range("P4").Value = 0.864385399390831
$ws.Range("Q4").Value = 147.4716204613254
$ws.Range("R4").Value = 1327.244584151928
$ws.Range("S4").Value = 0.01061748852979389
$ws.Range("T4").Value = 0.01061748852979389
$ws.Range("G5").Value = 7.974945333333333
$ws.Range("H5").Value = 23.924836
$ws.Range("I5").Value = 0.01228328074175765
$ws.Range("J5").Value = 0.01228328074175765
$ws.Range("M5").Value = 0.2265353333333333
$ws.Range("N5").Value = 0.6796059999999999
$ws.Range("O5").Value = 0.01058918740701822
$ws.Range("P5").Value = 0.01058918740701822
$ws.Range("Q5").Value = 1.806606899401778
$ws.Range("R5").Value = 16.259462094616
$ws.Range("S5").Value = 0.0001300699617474895
$ws.Range("T5").Value = 0.0001300699617474894
$ws.Range("I6").Value = 0.881245147037608
$ws.Range("J6").Value = 0.8812451470376081
$ws.Range("M6").Value = 0.3360566666666667
$ws.Range("N6").Value = 1.00817
$ws.Range("O6").Value = 0.01570866217798777
$ws.Range("P6").Value = 0.01570866217798777
$ws.Range("Q6").Value = 192.2748983422345
$ws.Range("R6").Value = 1730.47408508011
$ws.Range("S6").Value = 0.01384318231080495
$ws.Range("T6").Value = 0.01384318231080495
$ws.Range("I7").Value = 0.881245147037608
$ws.Range("J7").Value = 0.8812451470376081
$ws.Range("O7").Value = 0.109316751024163
$ws.Range("P7").Value = 0.1093167510241629
$ws.Range("S7").Value = 0.09633485632996207
$ws.Range("T7").Value = 0.09633485632996207
$ws.Range("I8").Value = 0.881245147037608
$ws.Range("J8").Value = 0.8812451470376081
$ws.Range("M8").Value = 18.491866
$ws.Range("N8").Value = 55.47559800000001
$ws.Range("O8").Value = 0.864385399390831
$ws.Range("P8").Value = 0.864385399390831
$ws.Range("Q8").Value = 10580.12534188149
$ws.Range("R8").Value = 95221.12807693344
$ws.Range("S8").Value = 0.7617354383833345
$ws.Range("T8").Value = 0.7617354383833346
$ws.Range("I9").Value = 0.881245147037608
$ws.Range("J9").Value = 0.8812451470376081
$ws.Range("M9").Value = 0.2265353333333333
$ws.Range("N9").Value = 0.6796059999999999
$ws.Range("O9").Value = 0.01058918740701822
$ws.Range("P9").Value = 0.01058918740701822
$ws.Range("Q9").Value = 129.6122425412109
$ws.Range("R9").Value = 1166.510182870898
$ws.Range("S9").Value = 0.009331670013506557
$ws.Range("T9").Value = 0.009331670013506557
$ws.Range("G10").Value = 68.85358966666666
$ws.Range("H10").Value = 206.560769
$ws.Range("I10").Value = 0.1060506293903268
$ws.Range("J10").Value = 0.1060506293903268
$ws.Range("M10").Value = 0.3360566666666667
$ws.Range("N10").Value = 1.00817
$ws.Range("O10").Value = 0.01570866217798777
$ws.Range("P10").Value = 0.01570866217798777
$ws.Range("Q10").Value = 23.13870783141444
$ws.Range("R10").Value = 208.24837048273
$ws.Range("S10").Value = 0.001665913510855626
$ws.Range("T10").Value = 0.001665913510855626
$ws.Range("G11").Value = 68.85358966666666
$ws.Range("H11").Value = 206.560769
$ws.Range("I11").Value = 0.1060506293903268
$ws.Range("J11").Value = 0.1060506293903268
$ws.Range("O11").Value = 0.109316751024163
$ws.Range("P11").Value = 0.1093167510241629
$ws.Range("Q11").Value = 161.0225195734393
$ws.Range("R11").Value = 1449.202676160954
$ws.Range("S11").Value = 0.01159311024901814
$ws.Range("T11").Value = 0.01159311024901814
$ws.Range("G12").Value = 68.85358966666666
$ws.Range("H12").Value = 206.560769
$ws.Range("I12").Value = 0.1060506293903268
$ws.Range("J12").Value = 0.1060506293903268
$ws.Range("M12").Value = 18.491866
$ws.Range("N12").Value = 55.47559800000001
$ws.Range("O12").Value = 0.864385399390831
$ws.Range("P12").Value = 0.864385399390831
$ws.Range("Q12").Value = 1273.231353734985
$ws.Range("R12").Value = 11459.08218361486
$ws.Range("S12").Value = 0.09166861564120668
$ws.Range("T12").Value = 0.09166861564120668
$ws.Range("G13").Value = 68.85358966666666
$ws.Range("H13").Value = 206.560769
$ws.Range("I13").Value = 0.1060506293903268
$ws.Range("J13").Value = 0.1060506293903268
$ws.Range("M13").Value = 0.2265353333333333
$ws.Range("N13").Value = 0.6796059999999999
$ws.Range("O13").Value = 0.01058918740701822
$ws.Range("P13").Value = 0.01058918740701822
$ws.Range("Q13").Value = 15.59777088633489
$ws.Range("R13").Value = 140.379937977014
$ws.Range("S13").Value = 0.001122989989246405
$ws.Range("T13").Value = 0.001122989989246405
$ws.Range("G14").Value = 0.273298
$ws.Range("H14").Value = 0.8198939999999999
$ws.Range("I14").Value = 0.0004209428303074948
$ws.Range("J14").Value = 0.0004209428303074948
$ws.Range("M14").Value = 0.3360566666666667
$ws.Range("N14").Value = 1.00817
$ws.Range("O14").Value = 0.01570866217798777
$ws.Range("P14").Value = 0.01570866217798777
$ws.Range("Q14").Value = 0.09184361488666666
$ws.Range("R14").Value = 0.8265925339799999
$ws.Range("S14").Value = 0.000006612448717546469
$ws.Range("T14").Value = 0.000006612448717546468
$ws.Range("G15").Value = 0.273298
$ws.Range("H15").Value = 0.8198939999999999
$ws.Range("I15").Value = 0.0004209428303074948
$ws.Range("J15").Value = 0.0004209428303074948
$ws.Range("O15").Value = 0.109316751024163
$ws.Range("P15").Value = 0.1093167510241629
$ws.Range("Q15").Value = 0.6391407153559999
$ws.Range("R15").Value = 5.752266438203999
$ws.Range("S15").Value = 0.00004601610257613089
$ws.Range("T15").Value = 0.00004601610257613088
$ws.Range("G16").Value = 0.273298
$ws.Range("H16").Value = 0.8198939999999999
$ws.Range("I16").Value = 0.0004209428303074948
$ws.Range("J16").Value = 0.0004209428303074948
$ws.Range("M16").Value = 18.491866
$ws.Range("N16").Value = 55.47559800000001
$ws.Range("O16").Value = 0.864385399390831
$ws.Range("P16").Value = 0.864385399390831
$ws.Range("Q16").Value = 5.053789994068
$ws.Range("R16").Value = 45.484109946612
$ws.Range("S16").Value = 0.0003638568364960507
$ws.Range("T16").Value = 0.0003638568364960507
$ws.Range("G17").Value = 0.273298
$ws.Range("H17").Value = 0.8198939999999999
$ws.Range("I17").Value = 0.0004209428303074948
$ws.Range("J17").Value = 0.0004209428303074948
$ws.Range("M17").Value = 0.2265353333333333
$ws.Range("N17").Value = 0.6796059999999999
$ws.Range("O17").Value = 0.01058918740701822
$ws.Range("P17").Value = 0.01058918740701822
$ws.Range("Q17").Value = 0.06191165352933332
$ws.Range("R17").Value = 0.5572048817639998
$ws.Range("S17").Value = 0.000004457442517766731
$ws.Range("T17").Value = 0.000004457442517766731
